$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet4: add a new data row (row 3) between the existing rows and the
# "Total" row (row 24). The SUM() formulas in row 24 already cover this
# range so their cached results will be recalculated automatically.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Sheet4")

$ws4.Range("A3").Value = "Crosses:"
$ws4.Range("B3").Value = 15484016
$ws4.Range("C3").Value = "Trials:"
$ws4.Range("D3").Value = 75000000
$ws4.Range("E3").Value = "p"
$ws4.Range("F3").Value = "'="
$ws4.Range("F3").Style = "Normal"
$ws4.Range("G3").Value = 0.20645354666666599
$ws4.Range("H3").Value = "Move"
$ws4.Range("I3").Value = "Avg:"
$ws4.Range("J3").Value = 12.3872128

# ---------------------------------------------------------------------------
# Sheet5: add two new data rows (rows 4 and 5).
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Sheet5")

$ws5.Range("A4").Value = "Crosses:"
$ws5.Range("B4").Value = 34047452
$ws5.Range("C4").Value = "Trials:"
$ws5.Range("D4").Value = 75000000
$ws5.Range("E4").Value = "p"
$ws5.Range("F4").Value = "'="
$ws5.Range("F4").Style = "Normal"
$ws5.Range("G4").Value = 0.45396602666666602
$ws5.Range("H4").Value = "Move"
$ws5.Range("I4").Value = "Avg:"
$ws5.Range("J4").Value = 4.5396602666666599

$ws5.Range("A5").Value = "Crosses:"
$ws5.Range("B5").Value = 34053806
$ws5.Range("C5").Value = "Trials:"
$ws5.Range("D5").Value = 75000000
$ws5.Range("E5").Value = "p"
$ws5.Range("F5").Value = "'="
$ws5.Range("F5").Style = "Normal"
$ws5.Range("G5").Value = 0.45405074666666601
$ws5.Range("H5").Value = "Move"
$ws5.Range("I5").Value = "Avg:"
$ws5.Range("J5").Value = 4.5405074666666598

# Sheet5 gets an explicit portrait page setup (it had none before).
$ws5.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Move the active tab/selection from Sheet5 to Sheet4: select F24 on Sheet5
# first (so its stored selection updates), then activate Sheet4 and select
# F24 there too, leaving Sheet4 as the active sheet.
# ---------------------------------------------------------------------------
$ws5.Activate()
$ws5.Range("F24").Select()

$ws4.Activate()
$ws4.Range("F24").Select()
